$wb = $excel.ActiveWorkbook

$wsHome = $wb.Worksheets.Item("HomePage")
$wsLogin = $wb.Worksheets.Item("LoginPage")
$wsReg = $wb.Worksheets.Item("RegistrationPage")

# --- LoginPage: add a second (duplicate) hyperlinked email column and a
#     new big "Test123" title cell, for parallel / cross-browser test runs ---
$wsLogin.Range("D1").Value = "abctest439@gmail.com"
$wsLogin.Hyperlinks.Add($wsLogin.Range("D1"), "mailto:abctest439@gmail.com") | Out-Null
$wsLogin.Range("D1").Style = "Hyperlink"

$wsLogin.Range("E1").Value = "Test123"
$wsHome.Range("A1").Copy()
$wsLogin.Range("E1").PasteSpecial(-4122)

$wsLogin.Range("A1:E1").RowHeight = 24

$wsLogin.Columns.Item(4).AutoFit() | Out-Null
$wsLogin.Columns.Item(5).AutoFit() | Out-Null

# --- HomePage: B1 "Good" -> "Poor" ---
$wsHome.Range("B1").Value = "Poor"
$wsHome.Range("B1").Select() | Out-Null

# --- RegistrationPage: its own lingering selection stays put, it's simply
#     no longer the active tab ---
$wsReg.Range("D6").Select() | Out-Null

# --- LoginPage becomes the active / selected sheet & cell ---
$wsLogin.Activate() | Out-Null
$wsLogin.Range("E1").Select() | Out-Null

$wb.Save() | Out-Null
